$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric value into a cell even when the cell is
# formatted as Text (numFmtId 49) - columns L/M on this sheet. Assigning
# .Value directly to a Text-formatted cell stores a text string (accurate
# Excel behaviour), so briefly switch to General, write the number, then
# restore the original format code.
function Set-NumericValue($range, $value) {
    $originalFormat = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $originalFormat
}

# Row 239: new positive cases updated 670 -> 671
$ws.Range("C239").Value = 671

# Row 242: new positive cases updated 522 -> 521
$ws.Range("C242").Value = 521

# Row 319: new positive cases updated 68 -> 69
$ws.Range("C319").Value = 69

# Row 320: new positive cases updated 71 -> 72
$ws.Range("C320").Value = 72

# Row 321: new positive cases updated 131 -> 195, plus one new
# hospital death recorded (L321 0 -> 1)
$ws.Range("C321").Value = 195
Set-NumericValue $ws.Range("L321") 1

# Row 322: new positive cases updated 21 -> 95
$ws.Range("C322").Value = 95

# Row 323: figures for this day filled in (previously a blank placeholder
# row awaiting data)
$ws.Range("C323").Value = 21
$ws.Range("E323").Value = 11
$ws.Range("F323").Value = 8
$ws.Range("G323").Value = 99
Set-NumericValue $ws.Range("L323") 0
Set-NumericValue $ws.Range("M323") 0
